# Add slides for climate lecture
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new cells for row 11 (10_climate lecture materials)
$ws.Range("E11").Value = "10_climate"
$ws.Range("F11").Value = "10_climate.R"
$ws.Range("G11").Value = "09_climate.RData"

# Update the active selection to G12
$ws.Range("G12").Select()
